# Commit: "Adding the RES Hourly Production Forecast to the Portfolio"
#
# 1) The forecast window rolled forward 26 days: every Timestamp in column A
#    moves from 2024-08-29 to 2024-09-24 (same time-of-day grid).
# 2) The RES (renewable) hourly production forecast numbers for rows 29-58
#    (Power, Power_MW, Next_Power_MW, Average_Power_MW, Energy_MWh) are refreshed
#    with the new forecast run output.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: roll the Timestamp column (A2:A97) forward from 2024-08-29 to 2024-09-24 ---
for ($r = 2; $r -le 97; $r++) {
    $cur = $ws.Cells.Item($r, 1).Value2
    $new = $cur.Replace("2024-08-29", "2024-09-24")
    $ws.Cells.Item($r, 1).Value2 = $new
}

# --- Step 2: refresh the RES forecast figures for rows 29-58 (columns B-F) ---
$ws.Cells.Item(29, 2).Value2 = 0
$ws.Cells.Item(29, 3).Value2 = 0
$ws.Cells.Item(29, 5).Value2 = 0
$ws.Cells.Item(29, 6).Value2 = 0
$ws.Cells.Item(30, 2).Value2 = 0
$ws.Cells.Item(30, 3).Value2 = 0
$ws.Cells.Item(30, 4).Value2 = 0
$ws.Cells.Item(30, 5).Value2 = 0
$ws.Cells.Item(30, 6).Value2 = 0
$ws.Cells.Item(31, 2).Value2 = 3255.24415588379
$ws.Cells.Item(31, 3).Value2 = 0.0032552441558837
$ws.Cells.Item(31, 4).Value2 = 0
$ws.Cells.Item(31, 5).Value2 = 0.0016276220779418
$ws.Cells.Item(31, 6).Value2 = 0.0004069055194854
$ws.Cells.Item(32, 2).Value2 = 29269.15754191081
$ws.Cells.Item(32, 3).Value2 = 0.0292691575419107
$ws.Cells.Item(32, 4).Value2 = 0.0032552441558837
$ws.Cells.Item(32, 5).Value2 = 0.0162622008488971
$ws.Cells.Item(32, 6).Value2 = 0.0040655502122242
$ws.Cells.Item(33, 2).Value2 = 130964.7746988932
$ws.Cells.Item(33, 3).Value2 = 0.1309647746988932
$ws.Cells.Item(33, 4).Value2 = 0.0292691575419107
$ws.Cells.Item(33, 5).Value2 = 0.08011696612040189
$ws.Cells.Item(33, 6).Value2 = 0.0200292415301004
$ws.Cells.Item(34, 2).Value2 = 305373.0834147135
$ws.Cells.Item(34, 3).Value2 = 0.3053730834147134
$ws.Cells.Item(34, 4).Value2 = 0.1309647746988932
$ws.Cells.Item(34, 5).Value2 = 0.2181689290568032
$ws.Cells.Item(34, 6).Value2 = 0.0545422322642007
$ws.Cells.Item(35, 2).Value2 = 592367.5735677083
$ws.Cells.Item(35, 3).Value2 = 0.5923675735677083
$ws.Cells.Item(35, 4).Value2 = 0.3053730834147134
$ws.Cells.Item(35, 5).Value2 = 0.4488703284912108
$ws.Cells.Item(35, 6).Value2 = 0.1122175821228026
$ws.Cells.Item(36, 2).Value2 = 871247.2115885416
$ws.Cells.Item(36, 3).Value2 = 0.8712472115885417
$ws.Cells.Item(36, 4).Value2 = 0.5923675735677083
$ws.Cells.Item(36, 5).Value2 = 0.731807392578125
$ws.Cells.Item(36, 6).Value2 = 0.1829518481445312
$ws.Cells.Item(37, 2).Value2 = 1143692.049804688
$ws.Cells.Item(37, 3).Value2 = 1.143692049804688
$ws.Cells.Item(37, 4).Value2 = 0.8712472115885417
$ws.Cells.Item(37, 5).Value2 = 1.007469630696614
$ws.Cells.Item(37, 6).Value2 = 0.2518674076741536
$ws.Cells.Item(38, 2).Value2 = 1390311.949544271
$ws.Cells.Item(38, 3).Value2 = 1.390311949544271
$ws.Cells.Item(38, 4).Value2 = 1.143692049804688
$ws.Cells.Item(38, 5).Value2 = 1.267001999674479
$ws.Cells.Item(38, 6).Value2 = 0.3167504999186197
$ws.Cells.Item(39, 2).Value2 = 1650165.860677083
$ws.Cells.Item(39, 3).Value2 = 1.650165860677083
$ws.Cells.Item(39, 4).Value2 = 1.390311949544271
$ws.Cells.Item(39, 5).Value2 = 1.520238905110677
$ws.Cells.Item(39, 6).Value2 = 0.3800597262776693
$ws.Cells.Item(40, 2).Value2 = 1930147.520833333
$ws.Cells.Item(40, 3).Value2 = 1.930147520833334
$ws.Cells.Item(40, 4).Value2 = 1.650165860677083
$ws.Cells.Item(40, 5).Value2 = 1.790156690755208
$ws.Cells.Item(40, 6).Value2 = 0.4475391726888019
$ws.Cells.Item(41, 2).Value2 = 2173008.897135417
$ws.Cells.Item(41, 3).Value2 = 2.173008897135417
$ws.Cells.Item(41, 4).Value2 = 1.930147520833334
$ws.Cells.Item(41, 5).Value2 = 2.051578208984375
$ws.Cells.Item(41, 6).Value2 = 0.5128945522460937
$ws.Cells.Item(42, 2).Value2 = 2382480.729166667
$ws.Cells.Item(42, 3).Value2 = 2.382480729166666
$ws.Cells.Item(42, 4).Value2 = 2.173008897135417
$ws.Cells.Item(42, 5).Value2 = 2.277744813151042
$ws.Cells.Item(42, 6).Value2 = 0.5694362032877602
$ws.Cells.Item(43, 2).Value2 = 2575254.546223958
$ws.Cells.Item(43, 3).Value2 = 2.575254546223958
$ws.Cells.Item(43, 4).Value2 = 2.382480729166666
$ws.Cells.Item(43, 5).Value2 = 2.478867637695313
$ws.Cells.Item(43, 6).Value2 = 0.619716909423828
$ws.Cells.Item(44, 2).Value2 = 2753371.261067708
$ws.Cells.Item(44, 3).Value2 = 2.753371261067708
$ws.Cells.Item(44, 4).Value2 = 2.575254546223958
$ws.Cells.Item(44, 5).Value2 = 2.664312903645833
$ws.Cells.Item(44, 6).Value2 = 0.6660782259114582
$ws.Cells.Item(45, 2).Value2 = 2905623.309895833
$ws.Cells.Item(45, 3).Value2 = 2.905623309895833
$ws.Cells.Item(45, 4).Value2 = 2.753371261067708
$ws.Cells.Item(45, 5).Value2 = 2.829497285481771
$ws.Cells.Item(45, 6).Value2 = 0.7073743213704426
$ws.Cells.Item(46, 2).Value2 = 3013935.970052083
$ws.Cells.Item(46, 3).Value2 = 3.013935970052083
$ws.Cells.Item(46, 4).Value2 = 2.905623309895833
$ws.Cells.Item(46, 5).Value2 = 2.959779639973958
$ws.Cells.Item(46, 6).Value2 = 0.7399449099934894
$ws.Cells.Item(47, 2).Value2 = 3155303.184895834
$ws.Cells.Item(47, 3).Value2 = 3.155303184895834
$ws.Cells.Item(47, 4).Value2 = 3.013935970052083
$ws.Cells.Item(47, 5).Value2 = 3.084619577473958
$ws.Cells.Item(47, 6).Value2 = 0.7711548943684896
$ws.Cells.Item(48, 2).Value2 = 3263945.100260416
$ws.Cells.Item(48, 3).Value2 = 3.263945100260417
$ws.Cells.Item(48, 4).Value2 = 3.155303184895834
$ws.Cells.Item(48, 5).Value2 = 3.209624142578125
$ws.Cells.Item(48, 6).Value2 = 0.8024060356445313
$ws.Cells.Item(49, 2).Value2 = 3391673.826822916
$ws.Cells.Item(49, 3).Value2 = 3.391673826822916
$ws.Cells.Item(49, 4).Value2 = 3.263945100260417
$ws.Cells.Item(49, 5).Value2 = 3.327809463541666
$ws.Cells.Item(49, 6).Value2 = 0.8319523658854164
$ws.Cells.Item(50, 2).Value2 = 3472410.662760417
$ws.Cells.Item(50, 3).Value2 = 3.472410662760416
$ws.Cells.Item(50, 4).Value2 = 3.391673826822916
$ws.Cells.Item(50, 5).Value2 = 3.432042244791667
$ws.Cells.Item(50, 6).Value2 = 0.8580105611979166
$ws.Cells.Item(51, 2).Value2 = 3457321.852864583
$ws.Cells.Item(51, 3).Value2 = 3.457321852864583
$ws.Cells.Item(51, 4).Value2 = 3.472410662760416
$ws.Cells.Item(51, 5).Value2 = 3.4648662578125
$ws.Cells.Item(51, 6).Value2 = 0.866216564453125
$ws.Cells.Item(52, 2).Value2 = 3517494.397135416
$ws.Cells.Item(52, 3).Value2 = 3.517494397135416
$ws.Cells.Item(52, 4).Value2 = 3.457321852864583
$ws.Cells.Item(52, 5).Value2 = 3.487408125
$ws.Cells.Item(52, 6).Value2 = 0.87185203125
$ws.Cells.Item(53, 2).Value2 = 3581752.436197916
$ws.Cells.Item(53, 3).Value2 = 3.581752436197916
$ws.Cells.Item(53, 4).Value2 = 3.517494397135416
$ws.Cells.Item(53, 5).Value2 = 3.549623416666666
$ws.Cells.Item(53, 6).Value2 = 0.8874058541666666
$ws.Cells.Item(54, 2).Value2 = 3582533.778645833
$ws.Cells.Item(54, 3).Value2 = 3.582533778645833
$ws.Cells.Item(54, 4).Value2 = 3.581752436197916
$ws.Cells.Item(54, 5).Value2 = 3.582143107421875
$ws.Cells.Item(54, 6).Value2 = 0.8955357768554686
$ws.Cells.Item(55, 2).Value2 = 3570155.350260416
$ws.Cells.Item(55, 3).Value2 = 3.570155350260416
$ws.Cells.Item(55, 4).Value2 = 3.582533778645833
$ws.Cells.Item(55, 5).Value2 = 3.576344564453125
$ws.Cells.Item(55, 6).Value2 = 0.8940861411132812
$ws.Cells.Item(56, 2).Value2 = 3528979.671875
$ws.Cells.Item(56, 3).Value2 = 3.528979671875
$ws.Cells.Item(56, 4).Value2 = 3.570155350260416
$ws.Cells.Item(56, 5).Value2 = 3.549567511067708
$ws.Cells.Item(56, 6).Value2 = 0.887391877766927
$ws.Cells.Item(57, 2).Value2 = 996627.2265625
$ws.Cells.Item(57, 3).Value2 = 0.9966272265625
$ws.Cells.Item(57, 4).Value2 = 3.528979671875
$ws.Cells.Item(57, 5).Value2 = 2.26280344921875
$ws.Cells.Item(57, 6).Value2 = 0.5657008623046873
$ws.Cells.Item(58, 4).Value2 = 0.9966272265625
$ws.Cells.Item(58, 5).Value2 = 0.49831361328125
$ws.Cells.Item(58, 6).Value2 = 0.1245784033203125
